$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2720319339633761
$ws.Range("C2").Value = 0.05594294268207989
$ws.Range("D2").Value = 0.03229110889917308
$ws.Range("F2").Value = 0.4137833934051329
$ws.Range("G2").Value = 0.2620255378731429
$ws.Range("H2").Value = 0.4511555071653817
$ws.Range("I2").Value = 0.331977983111452
$ws.Range("K2").Value = 0.2826643853282746
$ws.Range("N2").Value = 1.050901539845839
$ws.Range("O2").Value = 1.346909471831694
$ws.Range("B3").Value = 0.2377593010000396
$ws.Range("C3").Value = 0.0517745992824814
$ws.Range("D3").Value = 0.02844010956418686
$ws.Range("F3").Value = 0.4151881916490048
$ws.Range("G3").Value = 0.2644241519481554
$ws.Range("H3").Value = 0.4550221083554504
$ws.Range("I3").Value = 0.336354328377972
$ws.Range("K3").Value = 0.2468440625797541
$ws.Range("N3").Value = 1.056841322574996
$ws.Range("O3").Value = 1.35990921690059
$ws.Range("B4").Value = 0.2166499071478256
$ws.Range("C4").Value = 0.04919847600119454
$ws.Range("D4").Value = 0.02606286690498649
$ws.Range("F4").Value = 0.41632312449957
$ws.Range("G4").Value = 0.2661003526162062
$ws.Range("H4").Value = 0.4575814952060959
$ws.Range("I4").Value = 0.3392276114768116
$ws.Range("K4").Value = 0.2247539977302893
$ws.Range("N4").Value = 1.060835036831321
$ws.Range("O4").Value = 1.368704070192805
$ws.Range("B5").Value = 0.2080317740062583
$ws.Range("C5").Value = 0.04814455341097812
$ws.Range("D5").Value = 0.02509098336015114
$ws.Range("F5").Value = 0.4168541473307101
$ws.Range("G5").Value = 0.266834539826398
$ws.Range("H5").Value = 0.4586710956598701
$ws.Range("I5").Value = 0.3404453095842932
$ws.Range("K5").Value = 0.2157286027971708
$ws.Range("N5").Value = 1.06254981621062
$ws.Range("O5").Value = 1.372492473508537
$ws.Range("B6").Value = 0.2065998001364449
$ws.Range("C6").Value = 0.04796930315842474
$ws.Range("D6").Value = 0.02492941505254009
$ws.Range("F6").Value = 0.4169464632697029
$ws.Range("G6").Value = 0.2669595376430749
$ws.Range("H6").Value = 0.4588548404655981
$ws.Range("I6").Value = 0.3406503345467566
$ws.Range("K6").Value = 0.2142285428465698
$ws.Range("N6").Value = 1.062839831290816
$ws.Range("O6").Value = 1.373133882134283
$ws.Range("B7").Value = 0.2165337434060746
$ws.Range("C7").Value = 0.0491842790556376
$ws.Range("D7").Value = 0.02604977237049866
$ws.Range("F7").Value = 0.4163300085417632
$ws.Range("G7").Value = 0.2661100471780316
$ws.Range("H7").Value = 0.45759600108385
$ws.Range("I7").Value = 0.3392438442483634
$ws.Range("K7").Value = 0.2246323723966128
$ws.Range("N7").Value = 1.060857809253839
$ws.Range("O7").Value = 1.368754334191756
$ws.Range("B8").Value = 0.2602287949497679
$ws.Range("C8").Value = 0.0545092191158858
$ws.Range("D8").Value = 0.0309659573969796
$ws.Range("F8").Value = 0.4142112560331412
$ws.Range("G8").Value = 0.2628103343243708
$ws.Range("H8").Value = 0.4524502808625783
$ws.Range("I8").Value = 0.3334482995855286
$ws.Range("K8").Value = 0.2703339144138397
$ws.Range("N8").Value = 1.052877708447454
$ws.Range("O8").Value = 1.351223057392033
$ws.Range("B9").Value = 0.3453662785729534
$ws.Range("C9").Value = 0.06481555059501432
$ws.Range("D9").Value = 0.04050353077053614
$ws.Range("F9").Value = 0.4122164674274629
$ws.Range("G9").Value = 0.2579555586544302
$ws.Range("H9").Value = 0.4438278529230502
$ws.Range("I9").Value = 0.3235607956518258
$ws.Range("K9").Value = 0.3591647332554544
$ws.Range("N9").Value = 1.039973528127909
$ws.Range("O9").Value = 1.323294848837307
$ws.Range("B10").Value = 0.4075535444003719
$ws.Range("C10").Value = 0.07230142134436335
$ws.Range("D10").Value = 0.04744565905095044
$ws.Range("F10").Value = 0.4120668538142525
$ws.Range("G10").Value = 0.255376476409765
$ws.Range("H10").Value = 0.438385628886337
$ws.Range("I10").Value = 0.3171973449859511
$ws.Range("K10").Value = 0.4239180937684068
$ws.Range("N10").Value = 1.032158501836072
$ws.Range("O10").Value = 1.306709070107814
$ws.Range("B11").Value = 0.4357594633601707
$ws.Range("C11").Value = 0.07568754749435413
$ws.Range("D11").Value = 0.05058921424051732
$ws.Range("F11").Value = 0.4122843834146721
$ws.Range("G11").Value = 0.2544182221353211
$ws.Range("H11").Value = 0.4361031395633859
$ws.Range("I11").Value = 0.3144980952069716
$ws.Range("K11").Value = 0.4532593477778448
$ws.Range("N11").Value = 1.02896336935104
$ws.Range("O11").Value = 1.300017967868001
$ws.Range("B12").Value = 0.4464277004712187
$ws.Range("C12").Value = 0.07696694728116427
$ws.Range("D12").Value = 0.05177746470302225
$ws.Range("F12").Value = 0.4124077971576057
$ws.Range("G12").Value = 0.2540863123265353
$ws.Range("H12").Value = 0.4352665646573897
$ws.Range("I12").Value = 0.3135040826834619
$ws.Range("K12").Value = 0.4643528980196834
$ws.Range("N12").Value = 1.027805094108253
$ws.Range("O12").Value = 1.297607017028952
$ws.Range("B13").Value = 0.4441306828166489
$ws.Range("C13").Value = 0.07669153390176575
$ws.Range("D13").Value = 0.05152165020737698
$ws.Range("F13").Value = 0.4123793929839934
$ws.Range("G13").Value = 0.2541564173534638
$ws.Range("H13").Value = 0.4354455021069157
$ws.Range("I13").Value = 0.3137169094505357
$ws.Range("K13").Value = 0.4619644869684976
$ws.Range("N13").Value = 1.02805225389897
$ws.Range("O13").Value = 1.298120795143333
$ws.Range("B14").Value = 0.4366374050719628
$ws.Range("C14").Value = 0.07579286215442949
$ws.Range("D14").Value = 0.05068701570961309
$ws.Range("F14").Value = 0.4122937143959788
$ws.Range("G14").Value = 0.2543902949651695
$ws.Range("H14").Value = 0.4360337579687297
$ws.Range("I14").Value = 0.3144157532226348
$ws.Range("K14").Value = 0.4541723727940052
$ws.Range("N14").Value = 1.028867042788882
$ws.Range("O14").Value = 1.29981715567699
$ws.Range("B15").Value = 0.4320458759492567
$ws.Range("C15").Value = 0.07524202567820737
$ws.Range("D15").Value = 0.05017549655187281
$ws.Range("F15").Value = 0.4122465775733701
$ws.Range("G15").Value = 0.2545375850473732
$ws.Range("H15").Value = 0.4363976952494539
$ws.Range("I15").Value = 0.3148474797759224
$ws.Range("K15").Value = 0.4493971965189303
$ws.Range("N15").Value = 1.02937284790179
$ws.Range("O15").Value = 1.30087222204277
$ws.Range("B16").Value = 0.405708479981115
$ws.Range("C16").Value = 0.07207973557152059
$ws.Range("D16").Value = 0.04723992373787667
$ws.Range("F16").Value = 0.412058382217019
$ws.Range("G16").Value = 0.2554434305421935
$ws.Range("H16").Value = 0.438538680289156
$ws.Range("I16").Value = 0.3173776832883384
$ws.Range("K16").Value = 0.421998187625519
$ws.Range("N16").Value = 1.032374545170107
$ws.Range("O16").Value = 1.307163534648623
$ws.Range("B17").Value = 0.3895294567278142
$ws.Range("C17").Value = 0.07013478489336933
$ws.Range("D17").Value = 0.04543529518215905
$ws.Range("F17").Value = 0.4120160576462482
$ws.Range("G17").Value = 0.2560542327609525
$ws.Range("H17").Value = 0.4399015696014104
$ws.Range("I17").Value = 0.3189799716926025
$ws.Range("K17").Value = 0.4051596955484911
$ws.Range("N17").Value = 1.03430810529963
$ws.Range("O17").Value = 1.311241773211165
$ws.Range("B18").Value = 0.3802159066862316
$ws.Range("C18").Value = 0.0690142954969275
$ws.Range("D18").Value = 0.04439596414300695
$ws.Range("F18").Value = 0.4120185964288439
$ws.Range("G18").Value = 0.2564257825866392
$ws.Range("H18").Value = 0.4407036546213874
$ws.Range("I18").Value = 0.3199199673120994
$ws.Range("K18").Value = 0.3954638257934562
$ws.Range("N18").Value = 1.035454125752722
$ws.Range("O18").Value = 1.313667834371046
$ws.Range("B19").Value = 0.3770611809865443
$ws.Range("C19").Value = 0.06863460964288493
$ws.Range("D19").Value = 0.0440438335879918
$ws.Range("F19").Value = 0.4120240739732921
$ws.Range("G19").Value = 0.2565550565165537
$ws.Range("H19").Value = 0.4409783513535785
$ws.Range("I19").Value = 0.3202413935794581
$ws.Range("K19").Value = 0.3921791365123113
$ws.Range("N19").Value = 1.035847972232645
$ws.Range("O19").Value = 1.314503058670326
$ws.Range("B20").Value = 0.3912525544126311
$ws.Range("C20").Value = 0.07034201567530829
$ws.Range("D20").Value = 0.04562754184408391
$ws.Range("F20").Value = 0.4120177811003387
$ws.Range("G20").Value = 0.2559871174940653
$ws.Range("H20").Value = 0.4397546056136576
$ws.Range("I20").Value = 0.3188075007047182
$ws.Range("K20").Value = 0.4069533065030839
$ws.Range("N20").Value = 1.034098768065505
$ws.Range("O20").Value = 1.310799319847732
$ws.Range("B21").Value = 0.4388387138493783
$ws.Range("C21").Value = 0.07605690194505144
$ws.Range("D21").Value = 0.0509322268633241
$ws.Range("F21").Value = 0.4123177666845308
$ws.Range("G21").Value = 0.254320758826573
$ws.Range("H21").Value = 0.4358602199464698
$ws.Range("I21").Value = 0.3142097221397364
$ws.Range("K21").Value = 0.4564615823892098
$ws.Range("N21").Value = 1.02862631858661
$ws.Range("O21").Value = 1.299315559728811
$ws.Range("B22").Value = 0.4698644623110795
$ws.Range("C22").Value = 0.07977524252773094
$ws.Range("D22").Value = 0.05438660723744704
$ws.Range("F22").Value = 0.4127530236150179
$ws.Range("G22").Value = 0.253412178033507
$ws.Range("H22").Value = 0.4334767736839211
$ws.Range("I22").Value = 0.3113688081412676
$ws.Range("K22").Value = 0.4887166474628089
$ws.Range("N22").Value = 1.025350774578605
$ws.Range("O22").Value = 1.292526190914018
$ws.Range("B23").Value = 0.4533125036717252
$ws.Range("C23").Value = 0.07779224798076712
$ws.Range("D23").Value = 0.05254410891828343
$ws.Range("F23").Value = 0.4124988414329493
$ws.Range("G23").Value = 0.2538805757370142
$ws.Range("H23").Value = 0.4347340721702508
$ws.Range("I23").Value = 0.3128700449023043
$ws.Range("K23").Value = 0.4715110447843074
$ws.Range("N23").Value = 1.027071486781793
$ws.Range("O23").Value = 1.296084283753984
$ws.Range("B24").Value = 0.3904735793151133
$ws.Range("C24").Value = 0.07024833383280793
$ws.Range("D24").Value = 0.04554063281108256
$ws.Range("F24").Value = 0.4120169182228182
$ws.Range("G24").Value = 0.2560173967925792
$ws.Range("H24").Value = 0.439820990274761
$ws.Range("I24").Value = 0.3188854162220007
$ws.Range("K24").Value = 0.406142462267411
$ws.Range("N24").Value = 1.0341933022966
$ws.Range("O24").Value = 1.310999099379586
$ws.Range("B25").Value = 0.3223960337371068
$ws.Range("C25").Value = 0.06204230796497257
$ws.Range("D25").Value = 0.03793463408625541
$ws.Range("F25").Value = 0.4125249470886274
$ws.Range("G25").Value = 0.2590956343323754
$ws.Range("H25").Value = 0.4460035086386753
$ws.Range("I25").Value = 0.3260774303423517
$ws.Range("K25").Value = 0.3352212107272408
$ws.Range("N25").Value = 1.043171398412277
$ws.Range("O25").Value = 1.330159520021269
